$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ABB")

# --- Update "Period Ending" column D figures (latest year restated) ---
$ws.Range("D8").Value  = 25196000    # Total Revenue
$ws.Range("D9").Value  = 17422000    # Cost of Revenue
$ws.Range("D10").Value = 7774000     # Gross Profit
$ws.Range("D12").Value = 1013000     # Research Development
$ws.Range("D14").Value = -258000     # Non Recurring
$ws.Range("D17").Value = 22966000    # Total Operating Expenses
$ws.Range("D18").Value = 2230000     # Operating Income or Loss
$ws.Range("D20").Value = 126000      # Total Other Income/Expenses Net
$ws.Range("D21").Value = 3457000     # Earnings Before Interest And Taxes
$ws.Range("D22").Value = 254000      # Interest Expense
$ws.Range("D23").Value = 2102000     # Income Before Tax
$ws.Range("D24").Value = 583000      # Income Tax Expense
$ws.Range("D26").Value = 1519000     # Income After Tax
$ws.Range("D27").Value = 1367000     # Net Income From Continuing Ops
$ws.Range("D29").Value = 850000      # Discontinued Operations
$ws.Range("D32").Value = -126000     # Other Items
$ws.Range("D33").Value = 2217000     # Net Income
$ws.Range("D35").Value = 2217000     # Net Income Applicable To Common Shares

$ws.Range("D81").Value = 2217000     # Net Income (Cash Flow Statement section)

# --- Replace now-unavailable prior-year (column J) figures with "NA" ---
$ws.Range("J21").Value  = "NA"
$ws.Range("J83").Value  = "NA"
$ws.Range("J94").Value  = "NA"
$ws.Range("J100").Value = "NA"
$ws.Range("J101").Value = "NA"
